# Productive-Browsing.pptx edit
#
# The authors' name credits on the title slide (slide 1) were each
# split across several runs (e.g. "Ali " / "Haisam" / " Muhammad " /
# "Rafid" / " | 1405013") because of spell-check "err" markers on the
# misspelled-looking proper nouns. The commit simplifies each
# paragraph back down to a single plain run with the full line of
# text (dropping the err="1" spell-check flags in the process).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

$cr = [char]13

# First assign a throw-away value so the subsequent assignment is not
# a same-text no-op (the run/paragraph structure is only rebuilt when
# the text actually changes), then assign the real, final text.
$sh.TextFrame.TextRange.Text = "_"
$sh.TextFrame.TextRange.Text = "Ali Haisam Muhammad Rafid | 1405013" + $cr + "Md. Toufikuzzaman | 1405015"
